$d = $word.ActiveDocument

# --- 1. Remove the existing (hidden) _GoBack bookmark that currently sits
#        right after "2020/7/27" so it doesn't interfere with the insertion
#        below, then we'll recreate it in its new location afterwards.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- 2. Insert a new paragraph after the "Update 2020/7/27" paragraph
#        containing "Text text text text", reproducing the same run-split
#        pattern (first glyph of each "word" plain, the remainder hinted
#        as eastAsia) seen in the first paragraph.
$insertAt = $d.Paragraphs(1).Range.End - 1
$r = $d.Range($insertAt, $insertAt)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r><w:t>T</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>ext</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>text</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>text</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>text</w:t></w:r>' +
  '</w:p>'
$r.InsertXML($newParaXml)

# --- 3. Re-add the _GoBack bookmark as an empty range at the end of the
#        newly inserted paragraph (just before its paragraph mark).
$newParaEnd = $d.Paragraphs(2).Range.End - 1
$rEnd = $d.Range($newParaEnd, $newParaEnd)
$d.Bookmarks.Add("_GoBack", $rEnd)

# --- 4. Create the footnotes.xml / endnotes.xml parts (with just the
#        standard separator / continuation-separator boilerplate) by
#        adding a throw-away footnote and immediately deleting it again.
$tmpFootnote = $d.Footnotes.Add($d.Content, "", "x")
$tmpFootnote.Delete()
